# Applies:
#  1. Rename property "populationsCSV" (row 7) -> "populationsFolder" (leftover previous name fix)
#  2. Remove the "compoundPropertiesFile" row (row 14), shifting rows below it up by one
#  3. Update the active cell selection to C17 (matches resulting workbook state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix leftover previous name: populationsCSV -> populationsFolder
$ws.Range("A7").Value = "populationsFolder"

# 2. Remove the compoundPropertiesFile row entirely (row 14), shifting rows below up
$ws.Rows(14).Delete()

# 3. Update selection to match post-edit state
$ws.Range("C17").Select()
